# A new work-order record was added as row 8. It is, cell for cell, a
# duplicate of the existing row 4 except for column T ("Tecnico Assegnato"),
# which is reassigned to a different technician. Replicate that by copying
# row 4 onto row 8 and then only touching column T.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$srcRow = 4
$dstRow = 8

$srcRange = $ws.Range("A" + $srcRow + ":CP" + $srcRow)
$dstRange = $ws.Range("A" + $dstRow + ":CP" + $dstRow)

# Copy the whole row of values/types across (keeps every cell's text vs.
# number vs. boolean typing identical to row 4, with no extra formatting).
$srcRange.Copy()
$dstRange.PasteSpecial(-4163) | Out-Null

# A handful of source cells hold an empty string rather than being truly
# blank; a values-only paste of an empty string does not materialize a
# cell, so create those explicitly. A leading apostrophe forces Excel to
# store them as text (matching the rest of the row) instead of as blank
# cells, and the style is then reset to match the unformatted look of the
# rest of the row.
$emptyCols = @("C","D","K","N","O","Q","V","Y","Z","AA","AB","AK","AQ","AR","AY","AZ","BB","BH","BI","BJ","BK","BN","BO","BS","BT","BX","BZ","CA","CB","CC","CF","CG","CH","CJ","CK","CL","CM","CN","CO","CP")
$emptyAddrs = ($emptyCols | ForEach-Object { $_ + $dstRow }) -join ","
$emptyRange = $ws.Range($emptyAddrs)

foreach ($area in $emptyRange.Areas) {
    $area.Value = "'"
}
foreach ($area in $emptyRange.Areas) {
    $area.Style = $ws.Range("A" + $srcRow).Style
}

# Now set the one cell that actually differs: the assigned technician.
$ws.Range("T" + $dstRow).Value = "Maramao Percheseimorto"

# Reflect the cell that ended up selected after the new row was entered.
$ws.Range("T9").Select()
